# Hortaliza, Vega Monumental Concepción - Acelga
# Insert two new weekly observation rows (new date 2022-08-11 / serial 44784)
# right before the existing row 207, pushing the rest of the price history
# down by two rows (old row 207 -> new row 209, ... old row 289 -> new row 291).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 207-208; everything currently at/after row 207
# (through the old last row 289) shifts down to 209..291.
$ws.Rows("207:208").Insert()

# --- New row 207 ("Primera" quality) ---
$ws.Cells.Item(207,1).Value  = 11
$ws.Cells.Item(207,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(207,3).Value  = "Bíobío"
$ws.Cells.Item(207,4).Value  = 44784
$ws.Cells.Item(207,5).Value  = 8
$ws.Cells.Item(207,6).Value  = 100112009
$ws.Cells.Item(207,7).Value  = "Acelga"
$ws.Cells.Item(207,8).Value  = "Sin especificar"
$ws.Cells.Item(207,9).Value  = "Primera"
$ws.Cells.Item(207,10).Value = 200
$ws.Cells.Item(207,11).Value = 700
$ws.Cells.Item(207,12).Value = 800
$ws.Cells.Item(207,13).Value = 750
$ws.Cells.Item(207,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(207,15).Value = "Región de Ñuble"
$ws.Cells.Item(207,16).Value = 750
$ws.Cells.Item(207,17).Value = 1
$ws.Cells.Item(207,18).Value = "Hortaliza"

# --- New row 208 ("Segunda" quality) ---
$ws.Cells.Item(208,1).Value  = 11
$ws.Cells.Item(208,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(208,3).Value  = "Bíobío"
$ws.Cells.Item(208,4).Value  = 44784
$ws.Cells.Item(208,5).Value  = 8
$ws.Cells.Item(208,6).Value  = 100112009
$ws.Cells.Item(208,7).Value  = "Acelga"
$ws.Cells.Item(208,8).Value  = "Sin especificar"
$ws.Cells.Item(208,9).Value  = "Segunda"
$ws.Cells.Item(208,10).Value = 100
$ws.Cells.Item(208,11).Value = 600
$ws.Cells.Item(208,12).Value = 600
$ws.Cells.Item(208,13).Value = 600
$ws.Cells.Item(208,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(208,15).Value = "Región de Ñuble"
$ws.Cells.Item(208,16).Value = 600
$ws.Cells.Item(208,17).Value = 1
$ws.Cells.Item(208,18).Value = "Hortaliza"

Write-Host "Inserted rows 207-208; sheet now spans through row 291."
